$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, shifting existing rows 180:266 down to 181:267
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new price observation
$ws.Range("A180").Value = 11
$ws.Range("B180").Value = 'Vega Monumental Concepción'
$ws.Range("C180").Value = 'Bíobío'
$ws.Range("D180").Value2 = 44818
$ws.Range("D180").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E180").Value = 8
$ws.Range("F180").Value = 100112045
$ws.Range("G180").Value = 'Zapallo'
$ws.Range("H180").Value = 'Paine'
$ws.Range("I180").Value = '1a (guarda)'
$ws.Range("J180").Value = 1100
$ws.Range("K180").Value = 300
$ws.Range("L180").Value = 350
$ws.Range("M180").Value = 323
$ws.Range("N180").Value = '$/kilo (volumen en unidades)'
$ws.Range("O180").Value = 'Región de O''Higgins'
$ws.Range("P180").Value = 323
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = 'Hortaliza'
